# 1. Remove the stray "_GoBack" bookmark that sat after the first ":- " run.
$d = $word.ActiveDocument
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete() | Out-Null
}

# 2. Insert the two new paragraphs (a blank spacer + the new "performance" note)
#    right after the existing trailing paragraph, reconstructing that trailing
#    empty paragraph afterwards so the document still ends on a blank line.
$lastPara = $d.Paragraphs.Last
$insertionRange = $lastPara.Range
$insertionRange.Collapse(1)
$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">When </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>performance</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> gets down, reliability also gets down.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:-</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/></w:rPr><w:t>When American singer Michael Jackson died on June 25, 2009, the Twitter server crashed after users were updating their status to include the words &quot;Michael Jackson&quot; at a rate of 100,000 tweets per hour</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>[47</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> URL]</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00E07517" w:rsidRPr="00311D8E" w:rsidRDefault="00E07517"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p>'
$insertionRange.InsertXML($newParagraphsXml) | Out-Null

# 3. Re-create the "_GoBack" bookmark at the same relative spot in the new
#    paragraph: right before the italic Michael-Jackson sentence.
$newContentPara = $d.Paragraphs.Item(4)
$searchRange = $newContentPara.Range.Duplicate
$searchRange.Find.Execute("When American singer", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bookmarkRange = $d.Range($searchRange.Start, $searchRange.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
